# This edit adds one new weekly price record for "Cilantro" (Femacal de La
# Calera) just before the existing row that used to be row 296, shifting
# that row and every following row down by one (the sheet grows from
# A1:R409 to A1:R410).
#
# The new record reuses the same market/category/quality/unit/origin
# metadata as the row it is inserted in front of, but carries its own
# date and volume/price figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 296; everything from the old row 296
# down to the old row 409 shifts down to rows 297-410.
$ws.Rows(296).Insert()

$ws.Range("A296").Value = 3
$ws.Range("B296").Value = "Femacal de La Calera"
$ws.Range("C296").Value = "Coquimbo"
$ws.Range("D296").Value = 44795
$ws.Range("D296").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E296").Value = 5
$ws.Range("F296").Value = 100112040
$ws.Range("G296").Value = "Cilantro"
$ws.Range("H296").Value = "Sin especificar"
$ws.Range("I296").Value = "Primera"
$ws.Range("J296").Value = 280
$ws.Range("K296").Value = 4000
$ws.Range("L296").Value = 4500
$ws.Range("M296").Value = 4214
$ws.Range("N296").Value = "$/docena de atados (3 kilos)"
$ws.Range("O296").Value = "Provincia de Quillota"
$ws.Range("P296").Value = 1405
$ws.Range("Q296").Value = 3
$ws.Range("R296").Value = "Hortaliza"
